# "Fruta / hortaliza, semanal"
#
# A new weekly price-record row is inserted for Zanahoria (Vega Modelo de
# Temuco) just before the current row 102, pushing every following record
# down by one row (old row 102 -> new row 103, ... old row 187 -> new row
# 188) and extending the sheet's used range from A1:R187 to A1:R188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 102:187 down to 103:188, leaving a blank row 102 behind
# (Excel copies the formatting of the row above, same as an interactive
# "Insert Sheet Rows").
$ws.Rows(102).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(102, 1).Value  = 10
$ws.Cells.Item(102, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value  = "La Araucanía"
$ws.Cells.Item(102, 4).Value  = 44447
$ws.Cells.Item(102, 5).Value  = 9
$ws.Cells.Item(102, 6).Value  = 100114013
$ws.Cells.Item(102, 7).Value  = "Zanahoria"
$ws.Cells.Item(102, 8).Value  = "Sin especificar"
$ws.Cells.Item(102, 9).Value  = "Primera"
$ws.Cells.Item(102, 10).Value = 65
$ws.Cells.Item(102, 11).Value = 5500
$ws.Cells.Item(102, 12).Value = 5500
$ws.Cells.Item(102, 13).Value = 5500
$ws.Cells.Item(102, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(102, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(102, 16).Value = 275
$ws.Cells.Item(102, 17).Value = 20
$ws.Cells.Item(102, 18).Value = "Hortaliza"
